$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.713.75"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.14%  "

$ws.Range("D3").Value = "'1.885.14"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.33%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "'0.7919"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.46%  "

$ws.Range("D6").Value = "'241.44"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.28%  "

$ws.Range("D7").Value = "'1.001"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.11%  "

$ws.Range("D8").Value = "'0.3162"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.57%  "

$ws.Range("D9").Value = "'25.32"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.58%  "

$ws.Range("D10").Value = "'0.06975"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.77%  "

$ws.Range("D11").Value = "'0.08039"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.05%  "

$ws.Range("D12").Value = "'0.7634"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.28%  "

$ws.Range("D13").Value = "'1.879.51"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.55%  "

$ws.Range("D14").Value = "'5.288"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.16%  "

$ws.Range("D15").Value = "'92.04"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.03%  "

$ws.Range("D16").Value = "'29.744.30"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.01%  "

$ws.Range("D17").Value = "'13.83"
$ws.Range("D17").Style = "Normal"

$ws.Range("D18").Value = "'5.937"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.45%  "

$ws.Range("D19").Value = "'243.03"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.84%  "

$ws.Range("D20").Value = "'0.000007668"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.58%  "

$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").Value = "'1.001"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.07%  "

$ws.Range("B22").Value = "Chainlink"
$ws.Range("C22").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D22").Value = "'8.130"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +15.63%  "

$ws.Range("D23").Value = "'2.107.14"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.36%  "

$ws.Range("E24").Value = "  +0.02%  "

$ws.Range("D25").Value = "'0.1679"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.09%  "

$ws.Range("D26").Value = "'9.287"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.02%  "

$ws.Range("D27").Value = "'164.91"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.68%  "

$ws.Range("E28").Value = "  -2.08%  "

$ws.Range("D29").Value = "'2.046"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.30%  "

$ws.Range("D30").Value = "'1.391"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.15%  "

$ws.Range("D31").Value = "'1.532"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.96%  "

$ws.Range("D32").Value = "'4.378"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.46%  "

$ws.Range("D33").Value = "'0.05660"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.26%  "

$ws.Range("D34").Value = "'4.042"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.35%  "

$ws.Range("D35").Value = "'1.259"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.51%  "

$ws.Range("E36").Value = "  -0.69%  "

$ws.Range("D37").Value = "'0.9993"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.09%  "

$ws.Range("D38").Value = "'2.643"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.81%  "

$ws.Range("D39").Value = "'0.01907"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.63%  "

$ws.Range("D40").Value = "'2.764"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.06%  "

$ws.Range("D41").Value = "'0.4392"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.38%  "

$ws.Range("D42").Value = "'72.24"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.50%  "

$ws.Range("D43").Value = "'5.816"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.43%  "

$ws.Range("D44").Value = "'1.001"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.08%  "

$ws.Range("D45").Value = "'0.8350"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.93%  "

$ws.Range("D46").Value = "'102.77"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.43%  "

$ws.Range("D47").Value = "'1.019.50"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.96%  "

$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "'1.857"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.58%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'9.853"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.89%  "

$ws.Range("D50").Value = "'7.411"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.93%  "

$ws.Range("D51").Value = "'2.038.17"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.09%  "
